$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.074304938316345
$ws.Range("B1").Value = 1.793997406959534
$ws.Range("C1").Value = 5.279239177703857
$ws.Range("D1").Value = 0.7954637408256531
$ws.Range("E1").Value = 0.4952170550823212
